# [ADDITIONAL SCRAPING] add a "Player Info" sheet ahead of the existing
# "ODI Batting" / "ODI Bowling" sheets, and replace the MATCH_CARD_LINK
# (full scorecard URL) column on both match sheets with a plain
# MATCH_CODE column holding just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet --------------------------------
# Worksheets.Add() inserts before the active sheet, so do this FIRST
# (before grabbing references to the other sheets) -- otherwise those
# references end up pointing at the new sheet's slot once everything
# shifts over.
$infoWs = $wb.Worksheets.Add()
$infoWs.Name = "Player Info"

$infoWs.Range("A1").Value = "ID"
$infoWs.Range("B1").Value = "NAME"
$infoWs.Range("C1").Value = "BATTING_HAND"
$infoWs.Range("D1").Value = "BOWL_STYLE"

$headerRange = $infoWs.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Keep the ID looking like the source data (a text code, not a number).
$infoWs.Range("A2").NumberFormat = "@"
$infoWs.Range("A2").Value = "4697"
$infoWs.Range("B2").Value = "Asitha Madusanka Fernando"
$infoWs.Range("C2").Value = "Right Handed"
$infoWs.Range("D2").Value = "Right Arm Medium Fast"

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE --------------------
$battingWs = $wb.Worksheets.Item("ODI Batting")
$battingWs.Range("D1").Value = "MATCH_CODE"
$battingWs.Range("D2:D6").NumberFormat = "@"
$battingWs.Range("D2").Value = "4059"
$battingWs.Range("D3").Value = "4451"
$battingWs.Range("D4").Value = "4470"
$battingWs.Range("D5").Value = "4471"
$battingWs.Range("D6").Value = "4675"

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE --------------------
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$bowlingWs.Range("B1").Value = "MATCH_CODE"
$bowlingWs.Range("B2:B5").NumberFormat = "@"
$bowlingWs.Range("B2").Value = "4059"
$bowlingWs.Range("B3").Value = "4451"
$bowlingWs.Range("B4").Value = "4470"
$bowlingWs.Range("B5").Value = "4675"
